$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 3.1
$ws.Range("Q3").Value = 1.93
$ws.Range("R3").Value = 1.93
$ws.Range("G6").Value = 1.21
$ws.Range("H6").Value = 5.6
$ws.Range("I6").Value = 10.75
$ws.Range("J6").Value = 1.57
$ws.Range("K6").Value = 2.77
$ws.Range("L6").Value = 8.75
$ws.Range("M6").Value = 1.02
$ws.Range("N6").Value = 9.75
$ws.Range("O6").Value = 1.13
$ws.Range("P6").Value = 5.1
$ws.Range("Q6").Value = 1.42
$ws.Range("R6").Value = 2.67
$ws.Range("S6").Value = 1.24
$ws.Range("T6").Value = 3.65
$ws.Range("U6").Value = 1.88
$ws.Range("V6").Value = 1.82
$ws.Range("W6").Value = 9.5
$ws.Range("Y6").Value = 9.25
$ws.Range("Z6").Value = 7.8
$ws.Range("AA6").Value = 10
$ws.Range("AB6").Value = 25
$ws.Range("AC6").Value = 9.75
$ws.Range("AD6").Value = 12.5
$ws.Range("AE6").Value = 23
$ws.Range("AF6").Value = 90
$ws.Range("AG6").Value = 35
$ws.Range("AH6").Value = 100
$ws.Range("AI6").Value = 35
$ws.Range("AJ6").Value = 400
$ws.Range("AK6").Value = 150
$ws.Range("AL6").Value = 100
$ws.Range("AM6").Value = 600
$ws.Range("AN6").Value = 3.25
$ws.Range("AO6").Value = 5.1
$ws.Range("AP6").Value = 14
$ws.Range("AQ6").Value = 11.5
$ws.Range("AR6").Value = 30
$ws.Range("AS6").Value = 150
$ws.Range("AT6").Value = 3.65
$ws.Range("AU6").Value = 8.75
$ws.Range("AV6").Value = 70
$ws.Range("AW6").Value = 11.5
$ws.Range("AX6").Value = 65
$ws.Range("AY6").Value = 50
$ws.Range("AZ6").Value = 500
$ws.Range("BA6").Value = 400
$ws.Range("BB6").Value = 500
$ws.Range("H7").Value = 3.65
$ws.Range("I7").Value = 1.6
$ws.Range("K7").Value = 2.18
$ws.Range("L7").Value = 2.18
$ws.Range("N7").Value = 7.3
$ws.Range("O7").Value = 1.3
$ws.Range("P7").Value = 3.25
$ws.Range("Q7").Value = 1.9
$ws.Range("R7").Value = 1.85
$ws.Range("S7").Value = 1.4
$ws.Range("T7").Value = 2.72
$ws.Range("U7").Value = 1.9
$ws.Range("V7").Value = 1.8
$ws.Range("Y7").Value = 16.5
$ws.Range("AC7").Value = 7.3
$ws.Range("AD7").Value = 7.2
$ws.Range("AG7").Value = 6.4
$ws.Range("AH7").Value = 7.2
$ws.Range("AI7").Value = 8.25
$ws.Range("AK7").Value = 13
$ws.Range("AL7").Value = 28
$ws.Range("AT7").Value = 2.72
$ws.Range("AX7").Value = 7.8
$ws.Range("AY7").Value = 18
$ws.Range("AZ7").Value = 26
$ws.Range("BA7").Value = 60
